# Apply the upstream/downstream process renaming edit across both sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Process")

# --- Sheet1 ---
$ws1.Range("H2").Value = "['proc_2', 'proc_3', 'proc_4', 'proc_5', 'res_5']"
$ws1.Range("I2").Value = "['proc_6', 'proc_7', 'proc_8', 'proc_9', 'proc_10']"

$ws1.Range("I3").Value = "['proc_1']"
$ws1.Range("I4").Value = "['proc_1']"
$ws1.Range("I5").Value = "['proc_1']"
$ws1.Range("I6").Value = "['proc_1']"

$ws1.Range("H7").Value = "['proc_1']"
$ws1.Range("H8").Value = "['proc_1']"
$ws1.Range("H9").Value = "['proc_1']"
$ws1.Range("H10").Value = "['proc_1']"

# --- Process sheet ---
$ws2.Range("I2").Value = "['proc_2', 'proc_3', 'proc_4', 'proc_5', 'res_5']"
$ws2.Range("J2").Value = "['proc_6', 'proc_7', 'proc_8', 'proc_9', 'proc_10']"

$ws2.Range("J3").Value = "['proc_1']"
$ws2.Range("J4").Value = "['proc_1']"
$ws2.Range("J5").Value = "['proc_1']"
$ws2.Range("J6").Value = "['proc_1']"

$ws2.Range("I7").Value = "['proc_1']"
$ws2.Range("I8").Value = "['proc_1']"
$ws2.Range("I9").Value = "['proc_1']"
$ws2.Range("I10").Value = "['proc_1']"
